$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") '68.187.78'
Set-TextValue $ws.Range("E2") '  +1.24%  '
Set-TextValue $ws.Range("D3") '3.349.50'
Set-TextValue $ws.Range("E3") '  +0.85%  '
Set-TextValue $ws.Range("E4") '  -0.02%  '
Set-TextValue $ws.Range("D5") '583.07'
Set-TextValue $ws.Range("E5") '  +0.23%  '
Set-TextValue $ws.Range("D6") '176.88'
Set-TextValue $ws.Range("E6") '  +0.53%  '
Set-TextValue $ws.Range("D7") '1.00'
Set-TextValue $ws.Range("E7") '  +0.03%  '
Set-TextValue $ws.Range("E8") '  +0.27%  '
Set-TextValue $ws.Range("D9") '0.183'
Set-TextValue $ws.Range("E9") '  +2.32%  '
Set-TextValue $ws.Range("E10") '  +0.87%  '
Set-TextValue $ws.Range("D11") '47.98'
Set-TextValue $ws.Range("E11") '  +5.59%  '
Set-TextValue $ws.Range("E12") '  +1.50%  '
Set-TextValue $ws.Range("D13") '685.18'
Set-TextValue $ws.Range("E13") '  +4.03%  '
Set-TextValue $ws.Range("D14") '3.885.22'
Set-TextValue $ws.Range("E14") '  +0.56%  '
Set-TextValue $ws.Range("D15") '8.42'
Set-TextValue $ws.Range("E15") '  +0.12%  '
Set-TextValue $ws.Range("D16") '68.164.45'
Set-TextValue $ws.Range("E16") '  +0.91%  '
Set-TextValue $ws.Range("D17") '0.120'
Set-TextValue $ws.Range("E17") '  +1.22%  '
Set-TextValue $ws.Range("D18") '3.343.07'
Set-TextValue $ws.Range("E18") '  +0.61%  '
Set-TextValue $ws.Range("D19") '17.44'
Set-TextValue $ws.Range("E19") '  +0.36%  '
Set-TextValue $ws.Range("E20") '  +2.34%  '
Set-TextValue $ws.Range("E21") '  +0.56%  '
Set-TextValue $ws.Range("E22") '  +0.23%  '
Set-TextValue $ws.Range("D23") '16.94'
Set-TextValue $ws.Range("E23") '  -0.54%  '
Set-TextValue $ws.Range("D24") '100.41'
Set-TextValue $ws.Range("E24") '  +0.91%  '
Set-TextValue $ws.Range("D25") '3.91'
Set-TextValue $ws.Range("E25") '  +1.52%  '
Set-TextValue $ws.Range("E26") '  +1.19%  '
Set-TextValue $ws.Range("D27") '9.50'
Set-TextValue $ws.Range("D28") '32.99'
Set-TextValue $ws.Range("E28") '  -1.96%  '
Set-TextValue $ws.Range("D29") '8.51'
Set-TextValue $ws.Range("E29") '  +0.99%  '
Set-TextValue $ws.Range("D30") '6.93'
Set-TextValue $ws.Range("E30") '  -7.34%  '
Set-TextValue $ws.Range("D31") '561.77'
Set-TextValue $ws.Range("E31") '  -2.00%  '
Set-TextValue $ws.Range("D32") '11.07'
Set-TextValue $ws.Range("E32") '  +0.95%  '
Set-TextValue $ws.Range("E33") '  +1.09%  '
Set-TextValue $ws.Range("D34") '57.99'
Set-TextValue $ws.Range("E34") '  +2.53%  '
Set-TextValue $ws.Range("D36") '3.717.81'
Set-TextValue $ws.Range("E36") '  +0.66%  '
Set-TextValue $ws.Range("D37") '3.30'
Set-TextValue $ws.Range("E37") '  -1.90%  '
Set-TextValue $ws.Range("E38") '  +5.01%  '
Set-TextValue $ws.Range("D39") '34.78'
Set-TextValue $ws.Range("E39") '  +0.79%  '
Set-TextValue $ws.Range("E40") '  +1.62%  '
Set-TextValue $ws.Range("E42") '  +0.56%  '
Set-TextValue $ws.Range("E43") '  +0.62%  '
Set-TextValue $ws.Range("E44") '  -1.46%  '
Set-TextValue $ws.Range("D45") '0.0410'
Set-TextValue $ws.Range("E45") '  +1.17%  '
Set-TextValue $ws.Range("E46") '  +2.23%  '
Set-TextValue $ws.Range("E47") '  +0.36%  '
Set-TextValue $ws.Range("E48") '  -0.06%  '
Set-TextValue $ws.Range("E49") '  -0.21%  '
Set-TextValue $ws.Range("D50") '131.71'
Set-TextValue $ws.Range("E50") '  +3.21%  '
Set-TextValue $ws.Range("D51") '2.56'
Set-TextValue $ws.Range("E51") '  -1.70%  '
